# Actualización de clases grabadas
# Applies the content edits described by the diff:
#  1. Language (en-US) tagging on "UserStories_MVP_PockerPlanning" block
#  2. Language (en-US) tagging + bookmark cleanup on "Sprint Planning" block
#  3. New "Testing" classes content appended at the end of the document
#  4. Removal of the leftover "_GoBack" bookmark

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Tema: UserStories_MVP_PockerPlanning" + "Parte1: ..." paragraphs
#    get tagged as English (en-US) text.
# ---------------------------------------------------------------------
$pTema1  = $d.Paragraphs(12)
$pParte1 = $d.Paragraphs(13)
$pTema1.Range.LanguageID  = "en-US"
$pParte1.Range.LanguageID = "en-US"

# ---------------------------------------------------------------------
# 2) "Tema: Sprint Planning" + "Parte1: <link>" paragraphs get tagged
#    as English (en-US) text as well, and the stray _GoBack bookmark
#    (left over from the old cursor position) is removed.
# ---------------------------------------------------------------------
$pTema2  = $d.Paragraphs(26)
$pParte2 = $d.Paragraphs(27)
$pTema2.Range.LanguageID  = "en-US"
$pParte2.Range.LanguageID = "en-US"

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 3) Append the new "Testing" classes content at the end of the
#    document, reusing the last (empty) paragraph for the first line.
# ---------------------------------------------------------------------
$newParas = @(
    "Tema: Testing de Caja Negra",
    "",
    "Parte 1: https://youtu.be/OUY0N9cuz18",
    "",
    "Parte 2: https://youtu.be/QJkThF0MpDs",
    "",
    "Tema: Testing de Caja Blanca",
    "",
    "Parte 1: https://youtu.be/5kfUgtONLE0",
    "",
    "Parte 2: https://youtu.be/o0IJSCVb7t8",
    "",
    "Tema: Ejecución de Casos de Prueba",
    "",
    "Parte 1: https://youtu.be/4sxrbciyBZY",
    "",
    "Tema: Repaso para el parcial",
    "",
    "Parte1: https://youtu.be/jK4MykdHH40"
)

$lastIndex = $d.Paragraphs.Count
$cur = $d.Paragraphs($lastIndex)

$first = $true
foreach ($line in $newParas) {
    if (-not $first) {
        $cur.Range.InsertParagraphAfter()
        $lastIndex = $d.Paragraphs.Count
        $cur = $d.Paragraphs($lastIndex)
    }
    $first = $false

    if ($line -ne "") {
        $cur.Range.InsertAfter($line)
        $cur.Range.LanguageID = "en-US"
    }
    $cur.SpaceAfter = 0
}

Write-Host "done"
